$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.856.53'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = '1.641.64'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('E4').Value = '  -0.65%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.99'
$ws.Range('E5').Value = '  -0.88%  '
$ws.Range('E6').Value = '  +1.76%  '
$ws.Range('E7').Value = '  -0.71%  '
$ws.Range('E8').Value = '  +1.80%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0621'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.75'
$ws.Range('E10').Value = '  +4.10%  '
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('D12').Value = '1.870.99'
$ws.Range('E12').Value = '  +0.53%  '
$ws.Range('D13').Value = '1.640.71'
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.13'
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.528'
$ws.Range('E15').Value = '  +1.29%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '66.40'
$ws.Range('E16').Value = '  +4.03%  '
$ws.Range('D17').Value = '26.852.69'
$ws.Range('E17').Value = '  +0.66%  '
$ws.Range('E18').Value = '  +1.52%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '218.43'
$ws.Range('E19').Value = '  +3.93%  '
$ws.Range('E20').Value = '  -0.63%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.65'
$ws.Range('E21').Value = '  +8.11%  '
$ws.Range('E22').Value = '  +1.75%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.41'
$ws.Range('E23').Value = '  +3.90%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '145.91'
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('E26').Value = '  -0.83%  '
$ws.Range('E27').Value = '  +5.08%  '
$ws.Range('E28').Value = '  +1.25%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.80'
$ws.Range('E29').Value = '  +2.02%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0511'
$ws.Range('E30').Value = '  +2.27%  '
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.36'
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('E33').Value = '  +2.30%  '
$ws.Range('E34').Value = '  +2.97%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Value = '1.236.53'
$ws.Range('E36').Value = '  -1.87%  '
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('E38').Value = '  +3.94%  '
$ws.Range('E39').Value = '  +4.38%  '
$ws.Range('E40').Value = '  -0.68%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.805'
$ws.Range('E41').Value = '  +0.94%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.37'
$ws.Range('E42').Value = '  +2.49%  '
$ws.Range('D43').Value = '1.783.82'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('E44').Value = '  -2.84%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '60.78'
$ws.Range('E45').Value = '  +2.17%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '91.45'
$ws.Range('E46').Value = '  +0.61%  '
$ws.Range('E47').Value = '  +1.03%  '
$ws.Range('D48').Value = '0.0₆0104'
$ws.Range('E48').Value = '  +15.20%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0514'
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('E50').Value = '  +2.09%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.57'
$ws.Range('E51').Value = '  +2.28%  '
